$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition) - column F (想去人数), rows 2-10 all reset to 0
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 0
$ws1.Range("F3").Value = 0
$ws1.Range("F4").Value = 0
$ws1.Range("F5").Value = 0
$ws1.Range("F6").Value = 0
$ws1.Range("F7").Value = 0
$ws1.Range("F8").Value = 0
$ws1.Range("F9").Value = 0
$ws1.Range("F10").Value = 0

# Sheet 2: 演出 (Performance) - column F, rows 2-3 reset to 0
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F2").Value = 0
$ws2.Range("F3").Value = 0

# Sheet 3: 本地生活 (Local life) - header only, no data rows to update

# Sheet 4: 全部类型 (All types) - column F, rows 2-12 (row 5 unchanged, row 10 becomes 69)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 0
$ws4.Range("F3").Value = 0
$ws4.Range("F4").Value = 0
$ws4.Range("F6").Value = 0
$ws4.Range("F7").Value = 0
$ws4.Range("F8").Value = 0
$ws4.Range("F9").Value = 0
$ws4.Range("F10").Value = 69
$ws4.Range("F11").Value = 0
$ws4.Range("F12").Value = 0
